$wb = $excel.ActiveWorkbook

# Duplicate Sheet1 ("Title/example/date" table) into a new Sheet3 placed after Sheet2
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet1.Copy($null, $sheet2)

$newSheet = $wb.Worksheets.Item($sheet2.Index + 1)
$newSheet.Name = "Sheet3"

# Make the new sheet the active one with its own selection
$newSheet.Activate()
$newSheet.Range("E36").Select()
